$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3849203
$ws.Range("H79").Value = 3849203
$ws.Range("H98").Value = 935.5862
$ws.Range("I98").Value = 1054.2916
$ws.Range("J98").Value = 365.8
$ws.Range("K98").Value = 1054.2916
$ws.Range("L98").Value = 365.8
$ws.Range("M98").Value = 443.7084
$ws.Range("N98").Value = -3361.8
$ws.Range("H111").Value = 1122
$ws.Range("I111").Value = 944.3333
$ws.Range("J111").Value = 1477.3334
$ws.Range("K111").Value = 2832.9999
$ws.Range("L111").Value = 4432.0002
$ws.Range("M111").Value = 234.0001000000002
$ws.Range("N111").Value = -10566.0002
$ws.Range("H122").Value = 935.5862
$ws.Range("I122").Value = 1054.2916
$ws.Range("J122").Value = 365.8
$ws.Range("K122").Value = 3162.8748
$ws.Range("L122").Value = 1097.4
$ws.Range("M122").Value = -712.8748000000001
$ws.Range("N122").Value = -5997.4
$ws.Range("H137").Value = 788.7273
$ws.Range("I137").Value = 745.3333
$ws.Range("K137").Value = 2235.9999
$ws.Range("M137").Value = 314.0001000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2915.85
$ws.Range("I32").Value = 2254.628
$ws.Range("K32").Value = 2254.628
$ws.Range("M32").Value = -1967.628
$ws.Range("H61").Value = 823.0213
$ws.Range("I61").Value = 772.4286
$ws.Range("J61").Value = 1248
$ws.Range("K61").Value = 772.4286
$ws.Range("L61").Value = 1248
$ws.Range("M61").Value = -560.4286
$ws.Range("N61").Value = -1672
$ws.Range("H132").Value = 2035.3
$ws.Range("I132").Value = 1859.2858
$ws.Range("K132").Value = 5577.857400000001
$ws.Range("M132").Value = -3047.857400000001
$ws.Range("H136").Value = 823.0213
$ws.Range("I136").Value = 772.4286
$ws.Range("J136").Value = 1248
$ws.Range("K136").Value = 2317.2858
$ws.Range("L136").Value = 3744
$ws.Range("M136").Value = 232.7142000000003
$ws.Range("N136").Value = -8844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2254.9092
$ws.Range("I86").Value = 2030.5
$ws.Range("J86").Value = 4499
$ws.Range("K86").Value = 2030.5
$ws.Range("L86").Value = 4499
$ws.Range("M86").Value = -907.5
$ws.Range("N86").Value = -6745
$ws.Range("H89").Value = 2254.9092
$ws.Range("I89").Value = 2030.5
$ws.Range("J89").Value = 4499
$ws.Range("K89").Value = 10152.5
$ws.Range("L89").Value = 22495
$ws.Range("M89").Value = -4536.5
$ws.Range("N89").Value = -33727
$ws.Range("H105").Value = 3551.5417
$ws.Range("I105").Value = 4128.75
$ws.Range("J105").Value = 2397.125
$ws.Range("K105").Value = 4128.75
$ws.Range("L105").Value = 2397.125
$ws.Range("M105").Value = -2381.75
$ws.Range("N105").Value = -5891.125
$ws.Range("H107").Value = 9019.588
$ws.Range("I107").Value = 961.0833
$ws.Range("J107").Value = 28360
$ws.Range("K107").Value = 961.0833
$ws.Range("L107").Value = 28360
$ws.Range("M107").Value = 958.9167
$ws.Range("N107").Value = -32200
$ws.Range("H132").Value = 192188.58
$ws.Range("J132").Value = 192188.58
$ws.Range("L132").Value = 192188.58
$ws.Range("N132").Value = -202308.58
$ws.Range("H134").Value = 17138.262
$ws.Range("I134").Value = 1523.1321
$ws.Range("J134").Value = 86105.086
$ws.Range("K134").Value = 4569.3963
$ws.Range("L134").Value = 258315.258
$ws.Range("M134").Value = -2034.3963
$ws.Range("N134").Value = -263385.258

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2095553.4
$ws.Range("I31").Value = 3002634.2
$ws.Range("J31").Value = 2289.5386
$ws.Range("K31").Value = 3002634.2
$ws.Range("L31").Value = 2289.5386
$ws.Range("M31").Value = -3002339.2
$ws.Range("N31").Value = -2879.5386
$ws.Range("H34").Value = 2095553.4
$ws.Range("I34").Value = 3002634.2
$ws.Range("J34").Value = 2289.5386
$ws.Range("K34").Value = 3002634.2
$ws.Range("L34").Value = 2289.5386
$ws.Range("M34").Value = -3002432.2
$ws.Range("N34").Value = -2693.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 3930.6924
$ws.Range("I124").Value = 2699.5
$ws.Range("J124").Value = 4154.5454
$ws.Range("K124").Value = 8098.5
$ws.Range("L124").Value = 12463.6362
$ws.Range("M124").Value = -3188.5
$ws.Range("N124").Value = -22283.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3388.4375
$ws.Range("I80").Value = 3265.3572
$ws.Range("J80").Value = 4250
$ws.Range("K80").Value = 3265.3572
$ws.Range("L80").Value = 4250
$ws.Range("M80").Value = -2267.3572
$ws.Range("N80").Value = -6246
$ws.Range("H82").Value = 31124.875
$ws.Range("J82").Value = 31124.875
$ws.Range("L82").Value = 31124.875
$ws.Range("N82").Value = -31890.875
$ws.Range("H83").Value = 3388.4375
$ws.Range("I83").Value = 3265.3572
$ws.Range("J83").Value = 4250
$ws.Range("K83").Value = 16326.786
$ws.Range("L83").Value = 21250
$ws.Range("M83").Value = -11334.786
$ws.Range("N83").Value = -31234
$ws.Range("H85").Value = 31124.875
$ws.Range("J85").Value = 31124.875
$ws.Range("L85").Value = 31124.875
$ws.Range("N85").Value = -33776.875
$ws.Range("H122").Value = 10819597
$ws.Range("I122").Value = 9979021
$ws.Range("J122").Value = 12500750
$ws.Range("K122").Value = 29937063
$ws.Range("L122").Value = 37502250
$ws.Range("M122").Value = -29934613
$ws.Range("N122").Value = -37507150
$ws.Range("H124").Value = 38826.668
$ws.Range("J124").Value = 38826.668
$ws.Range("L124").Value = 38826.668
$ws.Range("N124").Value = -48646.668
$ws.Range("H132").Value = 2091.413
$ws.Range("I132").Value = 2058.12
$ws.Range("J132").Value = 2131.0476
$ws.Range("K132").Value = 6174.36
$ws.Range("L132").Value = 6393.1428
$ws.Range("M132").Value = -3644.36
$ws.Range("N132").Value = -11453.1428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5960.56
$ws.Range("I16").Value = 9073
$ws.Range("J16").Value = 1999.2727
$ws.Range("K16").Value = 9073
$ws.Range("L16").Value = 1999.2727
$ws.Range("M16").Value = -8903
$ws.Range("N16").Value = -2339.2727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 829.02325
$ws.Range("I122").Value = 781.8
$ws.Range("J122").Value = 1035.625
$ws.Range("K122").Value = 2345.4
$ws.Range("L122").Value = 3106.875
$ws.Range("M122").Value = 104.6000000000004
$ws.Range("N122").Value = -8006.875
$ws.Range("H123").Value = 49820
$ws.Range("J123").Value = 49820
$ws.Range("L123").Value = 49820
$ws.Range("N123").Value = -59620
$ws.Range("H125").Value = 35000
$ws.Range("J125").Value = 35000
$ws.Range("L125").Value = 35000
$ws.Range("N125").Value = -44840
